$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of trading/sentiment data (row 3)
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial()
$ws.Range("A3").Value = 42605.648333333331

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 24
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "Bag"
